$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("C1").Value = "orgao_responsavel"
$ws.Range("A1").Value = "sigla_om"

# --- Column D updates (indicativo_om) ---
$ws.Range("D3").Value = "CITBRA"
$ws.Range("D9").Value = "CFGOIA"
$ws.Range("D5").Value = "CPFBRA"

# --- Borders: replace old hair/mixed borders with a uniform thin box border ---
$rng = $ws.Range("A1:D10")
$rng.Borders.LineStyle = 1
$rng.Borders.Weight = 2
